$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---- Sheet1: new meeting-notes rows ----
# Row1 header cells become bold (Date / date value)
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("B1").Font.Bold = $true

# Row2: Attendance
$ws1.Range("A2").Value = "Attendance:"
$ws1.Range("A2").Font.Bold = $true
$ws1.Range("B2").Value = " Vincent ,Tommy, Yvonne"

# Row3: Activities
$ws1.Range("A3").Value = "Activities"
$ws1.Range("A3").Font.Bold = $true

# Row4 & 5: activity detail rows
$ws1.Range("A4").Value = "Go Through Project Spec"
$ws1.Range("A5").Value = "Select Topic - Recommender engine: create application for making recommendations based on user preferences."

# Row6: decision row
$ws1.Range("A6").Value = "Decide Project Process and tasks"

# shift the old "To do list:" block down to rows 8-9
$ws1.Range("A3").Copy() | Out-Null
$ws1.Range("A8").Value = "To do list:"
$ws1.Range("A8").Font.Bold = $true
$ws1.Range("A9").Value = "Vincent"
$ws1.Range("B9").Value = "Research,Find topics,data collection,http://surpriselib.com/"

# column widths / formatting
$ws1.Columns.Item(1).ColumnWidth = 33.66
$ws1.Columns.Item(2).ColumnWidth = 30.66

# page setup for Sheet1 (adds printerSettings-backed pageSetup like the other sheets)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---- Sheet2: move selection, drop tab-selected state ----
$ws2.Activate()
$ws2.Range("D2:G2").Select()

# ---- Sheet3: move selection only ----
$ws3.Activate()
$ws3.Range("S17").Select()

# ---- Sheet1 becomes the active / selected sheet again ----
$ws1.Activate()
$ws1.Range("B10").Select()
